# Updated the statistics after publication of 2.37.0
# Appends one new observation (row 35) to the "Data" table and fixes a
# previously mis-keyed value in the last existing row (X34).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$lo = $ws.ListObjects.Item("Data")

# --- Fix a typo in the last existing row (Command line options 136 -> 135)
$ws.Range("X34").Value = 135

# --- Grow the table by one row (table ref / autofilter / dimension all
#     follow automatically), then clone the formatting of the previous
#     last row onto the freshly appended one before filling in values.
$newListRow = $lo.ListRows.Add()

$ws.Range("A34:AF34").Copy()
$ws.Range("A35:AF35").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- New data row (2022-03-07, LoC=44790 serial date)
$ws.Range("A35").Value = 44790
$ws.Range("B35").Value = 106
$ws.Range("C35").Value = 242
$ws.Range("D35").Value = 195
$ws.Range("E35").Value = 4749
$ws.Range("F35").Value = 5535
$ws.Range("G35").Value = 1788
$ws.Range("H35").Value = 315
$ws.Range("I35").Value = 243
$ws.Range("J35").Value = 97
$ws.Range("K35").Value = 48
$ws.Range("L35").Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"
$ws.Range("M35").Value = 1791
$ws.Range("N35").Value = 3716
$ws.Range("O35").Value = 62325
$ws.Range("P35").Value = 43671
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("S35").Value = 225
$ws.Range("T35").Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"
$ws.Range("U35").Value = 0
$ws.Range("V35").Value = 160
$ws.Range("W35").Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"
$ws.Range("X35").Value = 136
$ws.Range("Y35").Value = 140
$ws.Range("Z35").Value = 4
$ws.Range("AA35").Value = 0
$ws.Range("AB35").Value = 322
$ws.Range("AC35").Value = 618
$ws.Range("AD35").Value = 7
# AE35 (GH runs) intentionally left blank, matching the source data.
$ws.Range("AF35").Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"

# --- Mirror the manual selection left behind after entering the new row.
$ws.Activate()
$ws.Range("AE35").Select()

$wb.Application.Calculate()
